# Add support for iGEM collections:
#  - new "IGEMSource" worksheet, inserted right after "EuroscarfSource"
#  - the "repository_name" validation lists on the various *Source sheets
#    gain the new "igem" option

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "IGEMSource" worksheet right after "EuroscarfSource"
# ---------------------------------------------------------------------
$euroscarf = $wb.Worksheets.Item("EuroscarfSource")
$igem = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $euroscarf)
$igem.Name = "IGEMSource"

# Header row (A1:H1)
$headers = @("sequence_file_url", "repository_id", "repository_name", "input", "output", "type", "output_name", "id")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $igem.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Data validation on column C (repository_name): list incl. "igem"
$igemRange = $igem.Range("C2:C1048576")
$igemRange.Validation.Add(3, 1, 1, '"addgene,genbank,benchling,snapgene,euroscarf,igem"')
$igemRange.Validation.ShowDropDown = $false
$igemRange.Validation.ShowInput = $false
$igemRange.Validation.ShowError = $false
$igemRange.Validation.IgnoreBlank = $true

# ---------------------------------------------------------------------
# 2. Extend the existing "repository_name"/"type" validation lists with
#    the new "igem" option on the sheets that already offer a repository
#    source drop-down.
# ---------------------------------------------------------------------
$newFormula = '"addgene,genbank,benchling,snapgene,euroscarf,igem"'

$repositoryIdSource = $wb.Worksheets.Item("RepositoryIdSource")
$repositoryIdSource.Range("B2:B1048576").Validation.Formula1 = $newFormula

$benchlingUrlSource = $wb.Worksheets.Item("BenchlingUrlSource")
$benchlingUrlSource.Range("B2:B1048576").Validation.Formula1 = $newFormula

$snapGenePlasmidSource = $wb.Worksheets.Item("SnapGenePlasmidSource")
$snapGenePlasmidSource.Range("B2:B1048576").Validation.Formula1 = $newFormula

$euroscarfSource = $wb.Worksheets.Item("EuroscarfSource")
$euroscarfSource.Range("B2:B1048576").Validation.Formula1 = $newFormula

$addGeneIdSource = $wb.Worksheets.Item("AddGeneIdSource")
$addGeneIdSource.Range("D2:D1048576").Validation.Formula1 = $newFormula

# Restore the originally active sheet/selection (the edit itself doesn't
# change which tab is active in the saved workbook).
$wb.Worksheets.Item(1).Select()
